$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1564649.4   # H15: 1410472.9 -> 1564649.4
$ws.Cells.Item(15, 9).Value = 1564649.4   # I15: 1410472.9 -> 1564649.4
$ws.Cells.Item(15, 11).Value = 4693948.199999999   # K15: 4231418.699999999 -> 4693948.199999999
$ws.Cells.Item(15, 13).Value = -4693779.199999999   # M15: -4231249.699999999 -> -4693779.199999999
$ws.Cells.Item(64, 8).Value = 52950   # H64: 53312.5 -> 52950
$ws.Cells.Item(64, 9).Value = 2921.4285   # I64: 3470.8333 -> 2921.4285
$ws.Cells.Item(64, 10).Value = 169683.33   # J64: 128075 -> 169683.33
$ws.Cells.Item(64, 11).Value = 2921.4285   # K64: 3470.8333 -> 2921.4285
$ws.Cells.Item(64, 12).Value = 169683.33   # L64: 128075 -> 169683.33
$ws.Cells.Item(64, 13).Value = -2673.4285   # M64: -3222.8333 -> -2673.4285
$ws.Cells.Item(64, 14).Value = -170179.33   # N64: -128571 -> -170179.33
$ws.Cells.Item(67, 8).Value = 52950   # H67: 53312.5 -> 52950
$ws.Cells.Item(67, 9).Value = 2921.4285   # I67: 3470.8333 -> 2921.4285
$ws.Cells.Item(67, 10).Value = 169683.33   # J67: 128075 -> 169683.33
$ws.Cells.Item(67, 11).Value = 2921.4285   # K67: 3470.8333 -> 2921.4285
$ws.Cells.Item(67, 12).Value = 169683.33   # L67: 128075 -> 169683.33
$ws.Cells.Item(67, 13).Value = -2063.4285   # M67: -2612.8333 -> -2063.4285
$ws.Cells.Item(67, 14).Value = -171399.33   # N67: -129791 -> -171399.33
$ws.Cells.Item(132, 8).Value = 1366.1786   # H132: 811.9846 -> 1366.1786
$ws.Cells.Item(132, 9).Value = 1481.4348   # I132: 860.38184 -> 1481.4348
$ws.Cells.Item(132, 10).Value = 836   # J132: 545.8 -> 836
$ws.Cells.Item(132, 11).Value = 4444.3044   # K132: 2581.14552 -> 4444.3044
$ws.Cells.Item(132, 12).Value = 2508   # L132: 1637.4 -> 2508
$ws.Cells.Item(132, 13).Value = -1914.3044   # M132: -51.14552000000003 -> -1914.3044
$ws.Cells.Item(132, 14).Value = -7568   # N132: -6697.4 -> -7568
$ws.Cells.Item(135, 8).Value = 1304.9149   # H135: 943.07465 -> 1304.9149
$ws.Cells.Item(135, 9).Value = 507.33334   # I135: 344.2157 -> 507.33334
$ws.Cells.Item(135, 10).Value = 2712.4119   # J135: 2851.9375 -> 2712.4119
$ws.Cells.Item(135, 11).Value = 4566.00006   # K135: 3097.9413 -> 4566.00006
$ws.Cells.Item(135, 12).Value = 24411.7071   # L135: 25667.4375 -> 24411.7071
$ws.Cells.Item(135, 13).Value = -2031.00006   # M135: -562.9413000000004 -> -2031.00006
$ws.Cells.Item(135, 14).Value = -29481.7071   # N135: -30737.4375 -> -29481.7071
$ws.Cells.Item(137, 8).Value = 2725.672   # H137: 2874.418 -> 2725.672
$ws.Cells.Item(137, 9).Value = 3142   # I137: 2854.1353 -> 3142
$ws.Cells.Item(137, 10).Value = 2235   # J137: 2899.4333 -> 2235
$ws.Cells.Item(137, 11).Value = 9426   # K137: 8562.4059 -> 9426
$ws.Cells.Item(137, 12).Value = 6705   # L137: 8698.2999 -> 6705
$ws.Cells.Item(137, 13).Value = -6876   # M137: -6012.4059 -> -6876
$ws.Cells.Item(137, 14).Value = -11805   # N137: -13798.2999 -> -11805
$ws.Cells.Item(141, 8).Value = 1928.541   # H141: 2062.5715 -> 1928.541
$ws.Cells.Item(141, 9).Value = 1542.1555   # I141: 1664.878 -> 1542.1555
$ws.Cells.Item(141, 10).Value = 3015.25   # J141: 3149.6 -> 3015.25
$ws.Cells.Item(141, 11).Value = 4626.4665   # K141: 4994.634 -> 4626.4665
$ws.Cells.Item(141, 12).Value = 9045.75   # L141: 9448.799999999999 -> 9045.75
$ws.Cells.Item(141, 13).Value = 553.5334999999995   # M141: 185.366 -> 553.5334999999995
$ws.Cells.Item(141, 14).Value = -19405.75   # N141: -19808.8 -> -19405.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1097.8334   # H2: 1014.5 -> 1097.8334
$ws.Cells.Item(2, 9).Value = 1077.4   # I2: 1096.75 -> 1077.4
$ws.Cells.Item(2, 10).Value = 1200   # J2: 850 -> 1200
$ws.Cells.Item(2, 11).Value = 1077.4   # K2: 1096.75 -> 1077.4
$ws.Cells.Item(2, 12).Value = 1200   # L2: 850 -> 1200
$ws.Cells.Item(2, 13).Value = -964.4000000000001   # M2: -983.75 -> -964.4000000000001
$ws.Cells.Item(2, 14).Value = -1426   # N2: -1076 -> -1426
$ws.Cells.Item(45, 8).Value = 1346.6957   # H45: 1327.1666 -> 1346.6957
$ws.Cells.Item(45, 9).Value = 1294   # I45: 1240.5714 -> 1294
$ws.Cells.Item(45, 10).Value = 1900   # J45: 1933.3334 -> 1900
$ws.Cells.Item(45, 11).Value = 1294   # K45: 1240.5714 -> 1294
$ws.Cells.Item(45, 12).Value = 1900   # L45: 1933.3334 -> 1900
$ws.Cells.Item(45, 13).Value = -917   # M45: -863.5714 -> -917
$ws.Cells.Item(45, 14).Value = -2654   # N45: -2687.3334 -> -2654
$ws.Cells.Item(46, 8).Value = 2184   # H46: 1888 -> 2184
$ws.Cells.Item(46, 9).Value = 3400   # I46: 2000 -> 3400
$ws.Cells.Item(46, 10).Value = 1576   # J46: 1850.6666 -> 1576
$ws.Cells.Item(46, 11).Value = 3400   # K46: 2000 -> 3400
$ws.Cells.Item(46, 12).Value = 1576   # L46: 1850.6666 -> 1576
$ws.Cells.Item(46, 13).Value = -3081   # M46: -1681 -> -3081
$ws.Cells.Item(46, 14).Value = -2214   # N46: -2488.6666 -> -2214
$ws.Cells.Item(61, 8).Value = 735.88635   # H61: 710.3095 -> 735.88635
$ws.Cells.Item(61, 9).Value = 741.3721   # I61: 710.3095 -> 741.3721
$ws.Cells.Item(61, 10).Value = 500   # J61: 0 -> 500
$ws.Cells.Item(61, 11).Value = 741.3721   # K61: 710.3095 -> 741.3721
$ws.Cells.Item(61, 12).Value = 500   # L61: 0 -> 500
$ws.Cells.Item(61, 13).Value = -529.3721   # M61: -498.3095 -> -529.3721
$ws.Cells.Item(61, 14).Value = -924   # N61: None -> -924
$ws.Cells.Item(63, 8).Value = 2507   # H63: 2419.8667 -> 2507
$ws.Cells.Item(66, 8).Value = 2507   # H66: 2419.8667 -> 2507
$ws.Cells.Item(74, 8).Value = 1511.3334   # H74: 1331.1177 -> 1511.3334
$ws.Cells.Item(74, 9).Value = 522.8570999999999   # I74: 549.61536 -> 522.8570999999999
$ws.Cells.Item(74, 10).Value = 4971   # J74: 3871 -> 4971
$ws.Cells.Item(74, 11).Value = 522.8570999999999   # K74: 549.61536 -> 522.8570999999999
$ws.Cells.Item(74, 12).Value = 4971   # L74: 3871 -> 4971
$ws.Cells.Item(74, 13).Value = 351.1429000000001   # M74: 324.38464 -> 351.1429000000001
$ws.Cells.Item(74, 14).Value = -6719   # N74: -5619 -> -6719
$ws.Cells.Item(76, 8).Value = 0   # H76: 17336 -> 0
$ws.Cells.Item(76, 10).Value = 0   # J76: 17336 -> 0
$ws.Cells.Item(76, 12).Value = ""   # L76: clear (was 17336)
$ws.Cells.Item(76, 14).Value = 0   # N76: -18012 -> 0
$ws.Cells.Item(77, 8).Value = 1511.3334   # H77: 1331.1177 -> 1511.3334
$ws.Cells.Item(77, 9).Value = 522.8570999999999   # I77: 549.61536 -> 522.8570999999999
$ws.Cells.Item(77, 10).Value = 4971   # J77: 3871 -> 4971
$ws.Cells.Item(77, 11).Value = 2614.2855   # K77: 2748.0768 -> 2614.2855
$ws.Cells.Item(77, 12).Value = 24855   # L77: 19355 -> 24855
$ws.Cells.Item(77, 13).Value = 1753.7145   # M77: 1619.9232 -> 1753.7145
$ws.Cells.Item(77, 14).Value = -33591   # N77: -28091 -> -33591
$ws.Cells.Item(79, 8).Value = 0   # H79: 17336 -> 0
$ws.Cells.Item(79, 10).Value = 0   # J79: 17336 -> 0
$ws.Cells.Item(79, 12).Value = ""   # L79: clear (was 17336)
$ws.Cells.Item(79, 14).Value = 0   # N79: -19676 -> 0
$ws.Cells.Item(116, 8).Value = 1097.8334   # H116: 1014.5 -> 1097.8334
$ws.Cells.Item(116, 9).Value = 1077.4   # I116: 1096.75 -> 1077.4
$ws.Cells.Item(116, 10).Value = 1200   # J116: 850 -> 1200
$ws.Cells.Item(116, 11).Value = 1077.4   # K116: 1096.75 -> 1077.4
$ws.Cells.Item(116, 12).Value = 1200   # L116: 850 -> 1200
$ws.Cells.Item(116, 13).Value = 1216.6   # M116: 1197.25 -> 1216.6
$ws.Cells.Item(116, 14).Value = -5788   # N116: -5438 -> -5788
$ws.Cells.Item(132, 8).Value = 1795.5088   # H132: 2009.84 -> 1795.5088
$ws.Cells.Item(132, 9).Value = 1809.7   # I132: 2089.2942 -> 1809.7
$ws.Cells.Item(132, 10).Value = 1762.1177   # J132: 1841 -> 1762.1177
$ws.Cells.Item(132, 11).Value = 5429.1   # K132: 6267.882599999999 -> 5429.1
$ws.Cells.Item(132, 12).Value = 5286.3531   # L132: 5523 -> 5286.3531
$ws.Cells.Item(132, 13).Value = -2899.1   # M132: -3737.882599999999 -> -2899.1
$ws.Cells.Item(132, 14).Value = -10346.3531   # N132: -10583 -> -10346.3531
$ws.Cells.Item(136, 8).Value = 735.88635   # H136: 710.3095 -> 735.88635
$ws.Cells.Item(136, 9).Value = 741.3721   # I136: 710.3095 -> 741.3721
$ws.Cells.Item(136, 10).Value = 500   # J136: 0 -> 500
$ws.Cells.Item(136, 11).Value = 2224.1163   # K136: 2130.9285 -> 2224.1163
$ws.Cells.Item(136, 12).Value = 1500   # L136: 0 -> 1500
$ws.Cells.Item(136, 13).Value = 325.8836999999999   # M136: 419.0715 -> 325.8836999999999
$ws.Cells.Item(136, 14).Value = -6600   # N136: None -> -6600

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1097.8334   # H3: 1014.5 -> 1097.8334
$ws.Cells.Item(3, 9).Value = 1077.4   # I3: 1096.75 -> 1077.4
$ws.Cells.Item(3, 10).Value = 1200   # J3: 850 -> 1200
$ws.Cells.Item(3, 11).Value = 1077.4   # K3: 1096.75 -> 1077.4
$ws.Cells.Item(3, 12).Value = 1200   # L3: 850 -> 1200
$ws.Cells.Item(3, 13).Value = -963.4000000000001   # M3: -982.75 -> -963.4000000000001
$ws.Cells.Item(3, 14).Value = -1428   # N3: -1078 -> -1428
$ws.Cells.Item(22, 8).Value = 259.7143   # H22: 248.27272 -> 259.7143
$ws.Cells.Item(22, 9).Value = 259.7143   # I22: 230.1 -> 259.7143
$ws.Cells.Item(22, 10).Value = 0   # J22: 430 -> 0
$ws.Cells.Item(22, 11).Value = 259.7143   # K22: 230.1 -> 259.7143
$ws.Cells.Item(22, 12).Value = 0   # L22: 430 -> 0
$ws.Cells.Item(22, 13).Value = ""   # M22: clear (was -57.09999999999999)
$ws.Cells.Item(22, 14).Value = -86.71429999999998   # N22: -776 -> -86.71429999999998
$ws.Cells.Item(64, 8).Value = 859635.5600000001   # H64: 467.5 -> 859635.5600000001
$ws.Cells.Item(64, 9).Value = 2577496   # I64: 148.66667 -> 2577496
$ws.Cells.Item(64, 10).Value = 705.375   # J64: 658.8 -> 705.375
$ws.Cells.Item(64, 11).Value = 2577496   # K64: 148.66667 -> 2577496
$ws.Cells.Item(64, 12).Value = 705.375   # L64: 658.8 -> 705.375
$ws.Cells.Item(64, 13).Value = -2577271   # M64: 76.33332999999999 -> -2577271
$ws.Cells.Item(64, 14).Value = -1155.375   # N64: -1108.8 -> -1155.375
$ws.Cells.Item(67, 8).Value = 859635.5600000001   # H67: 467.5 -> 859635.5600000001
$ws.Cells.Item(67, 9).Value = 2577496   # I67: 148.66667 -> 2577496
$ws.Cells.Item(67, 10).Value = 705.375   # J67: 658.8 -> 705.375
$ws.Cells.Item(67, 11).Value = 2577496   # K67: 148.66667 -> 2577496
$ws.Cells.Item(67, 12).Value = 705.375   # L67: 658.8 -> 705.375
$ws.Cells.Item(67, 13).Value = -2576716   # M67: 631.3333299999999 -> -2576716
$ws.Cells.Item(67, 14).Value = -2265.375   # N67: -2218.8 -> -2265.375
$ws.Cells.Item(132, 8).Value = 44890   # H132: 45053.332 -> 44890
$ws.Cells.Item(132, 10).Value = 44890   # J132: 45053.332 -> 44890
$ws.Cells.Item(132, 12).Value = 44890   # L132: 45053.332 -> 44890
$ws.Cells.Item(132, 14).Value = -55010   # N132: -55173.332 -> -55010
$ws.Cells.Item(134, 8).Value = 600.0566   # H134: 657.2157 -> 600.0566
$ws.Cells.Item(134, 9).Value = 573.13464   # I134: 593.8333 -> 573.13464
$ws.Cells.Item(134, 10).Value = 2000   # J134: 1671.3334 -> 2000
$ws.Cells.Item(134, 11).Value = 1719.40392   # K134: 1781.4999 -> 1719.40392
$ws.Cells.Item(134, 12).Value = 6000   # L134: 5014.0002 -> 6000
$ws.Cells.Item(134, 13).Value = 815.59608   # M134: 753.5001 -> 815.59608
$ws.Cells.Item(134, 14).Value = -11070   # N134: -10084.0002 -> -11070
$ws.Cells.Item(135, 8).Value = 40780   # H135: 33640 -> 40780
$ws.Cells.Item(135, 10).Value = 40780   # J135: 33640 -> 40780
$ws.Cells.Item(135, 12).Value = 40780   # L135: 33640 -> 40780
$ws.Cells.Item(135, 14).Value = -50920   # N135: -43780 -> -50920

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 27280.38   # H31: 30203.355 -> 27280.38
$ws.Cells.Item(31, 9).Value = 33690.156   # I31: 37166.137 -> 33690.156
$ws.Cells.Item(31, 10).Value = 15885.223   # J31: 17583.312 -> 15885.223
$ws.Cells.Item(31, 11).Value = 33690.156   # K31: 37166.137 -> 33690.156
$ws.Cells.Item(31, 12).Value = 15885.223   # L31: 17583.312 -> 15885.223
$ws.Cells.Item(31, 13).Value = -33395.156   # M31: -36871.137 -> -33395.156
$ws.Cells.Item(31, 14).Value = -16475.223   # N31: -18173.312 -> -16475.223
$ws.Cells.Item(34, 8).Value = 27280.38   # H34: 30203.355 -> 27280.38
$ws.Cells.Item(34, 9).Value = 33690.156   # I34: 37166.137 -> 33690.156
$ws.Cells.Item(34, 10).Value = 15885.223   # J34: 17583.312 -> 15885.223
$ws.Cells.Item(34, 11).Value = 33690.156   # K34: 37166.137 -> 33690.156
$ws.Cells.Item(34, 12).Value = 15885.223   # L34: 17583.312 -> 15885.223
$ws.Cells.Item(34, 13).Value = -33488.156   # M34: -36964.137 -> -33488.156
$ws.Cells.Item(34, 14).Value = -16289.223   # N34: -17987.312 -> -16289.223
$ws.Cells.Item(132, 8).Value = 861.2727   # H132: 1058.8206 -> 861.2727
$ws.Cells.Item(132, 9).Value = 736.5278   # I132: 943.9545000000001 -> 736.5278
$ws.Cells.Item(132, 10).Value = 1097.6316   # J132: 1207.4706 -> 1097.6316
$ws.Cells.Item(132, 11).Value = 2209.5834   # K132: 2831.8635 -> 2209.5834
$ws.Cells.Item(132, 12).Value = 3292.8948   # L132: 3622.4118 -> 3292.8948
$ws.Cells.Item(132, 13).Value = 320.4166   # M132: -301.8635000000004 -> 320.4166
$ws.Cells.Item(132, 14).Value = -8352.8948   # N132: -8682.4118 -> -8352.8948

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1210.1351   # H5: 1320.4286 -> 1210.1351
$ws.Cells.Item(5, 9).Value = 303.52942   # I5: 309.41177 -> 303.52942
$ws.Cells.Item(5, 10).Value = 1980.75   # J5: 2275.2778 -> 1980.75
$ws.Cells.Item(5, 11).Value = 910.58826   # K5: 928.23531 -> 910.58826
$ws.Cells.Item(5, 12).Value = 5942.25   # L5: 6825.8334 -> 5942.25
$ws.Cells.Item(5, 13).Value = -798.58826   # M5: -816.23531 -> -798.58826
$ws.Cells.Item(5, 14).Value = -6166.25   # N5: -7049.8334 -> -6166.25
$ws.Cells.Item(34, 8).Value = 3133.3333   # H34: 1940.6666 -> 3133.3333
$ws.Cells.Item(34, 9).Value = 5000   # I34: 1948 -> 5000
$ws.Cells.Item(34, 10).Value = 2760   # J34: 1933.3334 -> 2760
$ws.Cells.Item(34, 11).Value = 15000   # K34: 5844 -> 15000
$ws.Cells.Item(34, 12).Value = 8280   # L34: 5800.0002 -> 8280
$ws.Cells.Item(34, 13).Value = -14916   # M34: -5760 -> -14916
$ws.Cells.Item(34, 14).Value = -8448   # N34: -5968.0002 -> -8448
$ws.Cells.Item(39, 8).Value = 2185.76   # H39: 2214.3333 -> 2185.76
$ws.Cells.Item(39, 10).Value = 2185.76   # J39: 2214.3333 -> 2185.76
$ws.Cells.Item(39, 12).Value = 6557.280000000001   # L39: 6642.999899999999 -> 6557.280000000001
$ws.Cells.Item(39, 14).Value = -7145.280000000001   # N39: -7230.999899999999 -> -7145.280000000001
$ws.Cells.Item(47, 8).Value = 1278.9333   # H47: 1423.5625 -> 1278.9333
$ws.Cells.Item(47, 9).Value = 135.125   # I47: 96.59999999999999 -> 135.125
$ws.Cells.Item(47, 10).Value = 2586.1428   # J47: 2026.7273 -> 2586.1428
$ws.Cells.Item(47, 11).Value = 405.375   # K47: 289.8 -> 405.375
$ws.Cells.Item(47, 12).Value = 7758.428400000001   # L47: 6080.1819 -> 7758.428400000001
$ws.Cells.Item(47, 13).Value = 25.625   # M47: 141.2 -> 25.625
$ws.Cells.Item(47, 14).Value = -8620.428400000001   # N47: -6942.1819 -> -8620.428400000001
$ws.Cells.Item(55, 8).Value = 2121.8948   # H55: 3090.9092 -> 2121.8948
$ws.Cells.Item(55, 9).Value = 914.5   # I55: 0 -> 914.5
$ws.Cells.Item(55, 10).Value = 3000   # J55: 3090.9092 -> 3000
$ws.Cells.Item(55, 11).Value = 2743.5   # K55: 0 -> 2743.5
$ws.Cells.Item(55, 12).Value = 9000   # L55: 9272.7276 -> 9000
$ws.Cells.Item(55, 13).Value = -2566.5   # M55: None -> -2566.5
$ws.Cells.Item(55, 14).Value = -9354   # N55: -9626.7276 -> -9354
$ws.Cells.Item(107, 8).Value = 9672.362999999999   # H107: 9653.909 -> 9672.362999999999
$ws.Cells.Item(107, 9).Value = 14749   # I107: 14720 -> 14749
$ws.Cells.Item(107, 11).Value = 44247   # K107: 44160 -> 44247
$ws.Cells.Item(107, 13).Value = -42327   # M107: -42240 -> -42327
$ws.Cells.Item(122, 8).Value = 1503.3182   # H122: 1063.5862 -> 1503.3182
$ws.Cells.Item(122, 9).Value = 1303.909   # I122: 808.43475 -> 1303.909
$ws.Cells.Item(122, 10).Value = 1702.7273   # J122: 2041.6666 -> 1702.7273
$ws.Cells.Item(122, 11).Value = 11735.181   # K122: 7275.91275 -> 11735.181
$ws.Cells.Item(122, 12).Value = 15324.5457   # L122: 18374.9994 -> 15324.5457
$ws.Cells.Item(122, 13).Value = -9285.181   # M122: -4825.91275 -> -9285.181
$ws.Cells.Item(122, 14).Value = -20224.5457   # N122: -23274.9994 -> -20224.5457
$ws.Cells.Item(125, 8).Value = 8714.286   # H125: 7000 -> 8714.286
$ws.Cells.Item(125, 10).Value = 8714.286   # J125: 7000 -> 8714.286
$ws.Cells.Item(125, 12).Value = 26142.858   # L125: 21000 -> 26142.858
$ws.Cells.Item(125, 14).Value = -35982.858   # N125: -30840 -> -35982.858
$ws.Cells.Item(135, 8).Value = 1210.1351   # H135: 1320.4286 -> 1210.1351
$ws.Cells.Item(135, 9).Value = 303.52942   # I135: 309.41177 -> 303.52942
$ws.Cells.Item(135, 10).Value = 1980.75   # J135: 2275.2778 -> 1980.75
$ws.Cells.Item(135, 11).Value = 2731.76478   # K135: 2784.70593 -> 2731.76478
$ws.Cells.Item(135, 12).Value = 17826.75   # L135: 20477.5002 -> 17826.75
$ws.Cells.Item(135, 13).Value = -196.76478   # M135: -249.7059300000001 -> -196.76478
$ws.Cells.Item(135, 14).Value = -22896.75   # N135: -25547.5002 -> -22896.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 672.5961   # H132: 901.4167 -> 672.5961
$ws.Cells.Item(132, 9).Value = 540.625   # I132: 723.55554 -> 540.625
$ws.Cells.Item(132, 10).Value = 1112.5   # J132: 1435 -> 1112.5
$ws.Cells.Item(132, 11).Value = 1621.875   # K132: 2170.66662 -> 1621.875
$ws.Cells.Item(132, 12).Value = 3337.5   # L132: 4305 -> 3337.5
$ws.Cells.Item(132, 13).Value = 908.125   # M132: 359.33338 -> 908.125
$ws.Cells.Item(132, 14).Value = -8397.5   # N132: -9365 -> -8397.5
$ws.Cells.Item(136, 8).Value = 453.1111   # H136: 393.41177 -> 453.1111
$ws.Cells.Item(136, 9).Value = 414.9565   # I136: 393.41177 -> 414.9565
$ws.Cells.Item(136, 10).Value = 672.5   # J136: 0 -> 672.5
$ws.Cells.Item(136, 11).Value = 1244.8695   # K136: 1180.23531 -> 1244.8695
$ws.Cells.Item(136, 12).Value = 2017.5   # L136: 0 -> 2017.5
$ws.Cells.Item(136, 13).Value = 1305.1305   # M136: 1369.76469 -> 1305.1305
$ws.Cells.Item(136, 14).Value = -7117.5   # N136: None -> -7117.5

Write-Output "Applied Yojimbo Profits update across ALC/ARM/BSM/CRP/CUL/WVR sheets"